# feat: add 2022-Q4 data
#
# A new "2022-Q4" fund-holdings sheet is inserted right after the "总计"
# summary sheet (i.e. before the pre-existing "2021-Q3" sheet, which in
# turn sits before "2020-Q4"). The "总计" summary sheet grows a row so it
# keeps one row per quarter sheet, newest quarter on top.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$q4Sheet = $wb.Worksheets.Item("2020-Q4")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Header row — same headings as the other quarter sheets, except the
# 3rd column is "基金规模" here instead of "基金金额".
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match the bold/bordered header formatting used on the other quarter sheets.
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data rows. Fund codes and the percentage-looking figures are stored as
# text on the source sheet (leading zeros / trailing decimals that must
# round-trip verbatim), so they're entered with a leading apostrophe to
# keep Excel from re-interpreting them as numbers.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010387"
$newSheet.Range("C2").Value = "易方达医药生物股票A"
$newSheet.Range("D2").Value = "'20.77"
$newSheet.Range("E2").Value = "'91.88"
$newSheet.Range("F2").Value = "'5.04"
$newSheet.Range("G2").Value = "'1.0468"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'010388"
$newSheet.Range("C3").Value = "易方达医药生物股票C"
$newSheet.Range("D3").Value = "'6.05"
$newSheet.Range("E3").Value = "'91.88"
$newSheet.Range("F3").Value = "'5.04"
$newSheet.Range("G3").Value = "'0.3049"
$newSheet.Range("H3").Value = 4

# Re-flatten the cell format on the just-entered text cells back to the
# sheet's plain (unstyled) body format — typing an apostrophe-escaped
# value tags the cell with a "stored as text" number format, which the
# surrounding data cells don't otherwise carry.
$newSheet.Range("C2:C3").Copy()
$newSheet.Range("B2:B3").PasteSpecial(-4122)
$newSheet.Range("D2:G3").PasteSpecial(-4122)

# A2:A3 / H2:H3 stay plain numbers — copy the index-column style from the
# "2021-Q3" sheet so the leading index column keeps its bold/bordered look.
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2. Grow the "总计" (summary) sheet with a 2022-Q4 row, pushing the
#    existing 2021-Q3 / 2020-Q4 rows down by one.
# ---------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.05

$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.13

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 1.35
